$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to start its data at row 2 (row 1 was blank). Delete that
# leading blank row so the header moves up to row 1.
$ws.Rows(1).Delete()

# Insert six new task rows right after the "research different bike
# components" / "criteria" rows (old rows 6-7, now rows 5-6), pushing the
# remainder of the backlog down.
$ws.Rows("6:11").Insert()

# Populate the newly inserted rows (dates stay blank, only the task text is
# filled in) using the same italic "sub-task" style already used elsewhere
# in column C (B1/C1 header cells).
$ws.Range("C6").Value = "Research Mechanical Design Synthesis"
$ws.Range("C7").Value = "Develop Design Specification for kit "
$ws.Range("C8").Value = "Carry out morphological design for kit"
$ws.Range("C9").Value = "Sketch final concept"
$ws.Range("C10").Value = "Develop criteria based on concept, specification, and last year's report"
$ws.Range("C11").Value = "Decide on weights and scoring system"
$ws.Range("C6:C11").Font.Italic = $true

# Turn on the AutoFilter over the full table (one blank row past the data,
# matching the workbook's own A1:C57 filter range).
$ws.Range("A1:C57").AutoFilter() | Out-Null

# AutoFilter also registers the hidden workbook-level "_FilterDatabase" name
# scoped to this sheet.
$flt = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$57")
$flt.Visible = $false

# Restore the author's last active selection.
$ws.Range("G11").Select()
